$wb = $excel.ActiveWorkbook

# --- 1. Sheets: rename existing "High Priority break-up" and add a copy ---
$wsHP = $wb.Worksheets.Item("High Priority break-up")
$wsHP.Name = "Interannual update - High Pri"

# Duplicate it (preserves exact old content/formatting) and place the copy
# right after, then rename the copy - this becomes the "Major update" tab
# holding the historical (pre-update) numbers.
$wsHP.Copy([System.Reflection.Missing]::Value, $wsHP)
$wsMajor = $wb.Worksheets.Item("Interannual update - High Pri (2)")
$wsMajor.Name = "Major update - High Priority "

# --- 2. Replace the data on the renamed sheet with the new break-up figures ---
$wsHP.Range("A2:E2").ClearContents()

$wsHP.Range("A2").Value = "Trend New"
$wsHP.Range("B2").Value = 74
$wsHP.Range("C2").Value = 71.8
$wsHP.Range("D2").Value = 74
$wsHP.Range("E2").Value = 77.90000000000001

$wsHP.Range("A3").Value = "IUCN"
$wsHP.Range("B3").Value = 29
$wsHP.Range("C3").Value = 28.2
$wsHP.Range("D3").Value = 21
$wsHP.Range("E3").Value = 22.1

# --- 3. Trends Status sheet updates ---
$wsTrends = $wb.Worksheets.Item("Trends Status")
$wsTrends.Range("D2").ClearContents()
$wsTrends.Range("D3").ClearContents()
$wsTrends.Range("B4").Value = 0
$wsTrends.Range("D4").ClearContents()
$wsTrends.Range("D5").ClearContents()
$wsTrends.Range("D6").ClearContents()
$wsTrends.Range("B7").Value = 1
$wsTrends.Range("B8").Value = 363
$wsTrends.Range("C8").Value = 361

# --- 4. Priority Status sheet updates ---
$wsPriority = $wb.Worksheets.Item("Priority Status")
$wsPriority.Range("B2").Value = 103
$wsPriority.Range("B3").Value = 286
$wsPriority.Range("B4").Value = 554

# --- 5. Species qualification sheet updates ---
$wsSpecies = $wb.Worksheets.Item("Species qualification")
$wsSpecies.Range("A2").Value = "SoIB Assessment"
$wsSpecies.Range("B2").Value = 364
$wsSpecies.Range("C3").Value = 0
